# Wheel Strategy detection logic
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B4").Value = 48
$wsSummary.Range("B5").Value = 2
$wsSummary.Range("B6").Value = "Green Flag"
$wsSummary.Range("B7").Value = 4

# --- Symbols sheet ---
$wsSymbols = $wb.Worksheets.Item("Symbols")
$wsSymbols.Range("B2").Value = 48

# --- Strategies sheet ---
$wsStrategies = $wb.Worksheets.Item("Strategies")
$wsStrategies.Range("D2").Value = 48
$wsStrategies.Range("E2").Value = 50
$wsStrategies.Range("F2").Value = 2
$wsStrategies.Range("H2").Value = 12
$wsStrategies.Range("J2").Value = "[{'strategy_name': 'Short Put', 'pnl': 50.0, 'entry_ts': '2025-01-01T00:00:00', 'exit_ts': '2025-01-05T00:00:00'}]"
